$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 62 (old rows 62-70 shift down to 63-71)
$ws.Rows.Item(62).Insert()

# Populate the new row 62 with the new test step data
$ws.Cells.Item(62, 1).Value = "UploadDocumentScreen"
$ws.Cells.Item(62, 2).Value = "date_RevisionDate"
$ws.Cells.Item(62, 3).Value = "//input[@id='revisiondate_0_da']"

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 48
$ws.Range("B62").Select()
